# Auto-generated Excel COM-interop script to apply numeric cell updates
# per the target diff, across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6582855
$ws.Range("I74").Value = 3372.5715
$ws.Range("J74").Value = 10420887
$ws.Range("K74").Value = 3372.5715
$ws.Range("L74").Value = 10420887
$ws.Range("M74").Value = -2436.5715
$ws.Range("N74").Value = -10422759
$ws.Range("H77").Value = 6582855
$ws.Range("I77").Value = 3372.5715
$ws.Range("J77").Value = 10420887
$ws.Range("K77").Value = 16862.8575
$ws.Range("L77").Value = 52104435
$ws.Range("M77").Value = -12182.8575
$ws.Range("N77").Value = -52113795
$ws.Range("H96").Value = 17857774
$ws.Range("I96").Value = 25000572
$ws.Range("J96").Value = 776
$ws.Range("K96").Value = 75001716
$ws.Range("L96").Value = 2328
$ws.Range("M96").Value = -75000343
$ws.Range("N96").Value = -5074
$ws.Range("H100").Value = 2670
$ws.Range("I100").Value = 1880
$ws.Range("J100").Value = 2933.3333
$ws.Range("K100").Value = 1880
$ws.Range("L100").Value = 2933.3333
$ws.Range("M100").Value = -1339
$ws.Range("N100").Value = -4015.3333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 583.51514
$ws.Range("I2").Value = 591.96875
$ws.Range("J2").Value = 313
$ws.Range("K2").Value = 591.96875
$ws.Range("L2").Value = 313
$ws.Range("M2").Value = -478.96875
$ws.Range("N2").Value = -539
$ws.Range("H63").Value = 4466448.5
$ws.Range("I63").Value = 2630
$ws.Range("J63").Value = 15625995
$ws.Range("K63").Value = 2630
$ws.Range("L63").Value = 15625995
$ws.Range("M63").Value = -1944
$ws.Range("N63").Value = -15627367
$ws.Range("H66").Value = 4466448.5
$ws.Range("I66").Value = 2630
$ws.Range("J66").Value = 15625995
$ws.Range("K66").Value = 13150
$ws.Range("L66").Value = 78129975
$ws.Range("M66").Value = -9718
$ws.Range("N66").Value = -78136839
$ws.Range("H74").Value = 32259900
$ws.Range("I74").Value = 58824224
$ws.Range("J74").Value = 3221.2856
$ws.Range("K74").Value = 58824224
$ws.Range("L74").Value = 3221.2856
$ws.Range("M74").Value = -58823350
$ws.Range("N74").Value = -4969.2856
$ws.Range("H77").Value = 32259900
$ws.Range("I77").Value = 58824224
$ws.Range("J77").Value = 3221.2856
$ws.Range("K77").Value = 294121120
$ws.Range("L77").Value = 16106.428
$ws.Range("M77").Value = -294116752
$ws.Range("N77").Value = -24842.428
$ws.Range("H97").Value = 66667920
$ws.Range("I97").Value = 1083.9166
$ws.Range("K97").Value = 1083.9166
$ws.Range("M97").Value = -587.9166
$ws.Range("H102").Value = 1890
$ws.Range("I102").Value = 1800
$ws.Range("J102").Value = 1980
$ws.Range("K102").Value = 1800
$ws.Range("L102").Value = 1980
$ws.Range("M102").Value = -178
$ws.Range("N102").Value = -5224
$ws.Range("H116").Value = 583.51514
$ws.Range("I116").Value = 591.96875
$ws.Range("J116").Value = 313
$ws.Range("K116").Value = 591.96875
$ws.Range("L116").Value = 313
$ws.Range("M116").Value = 1702.03125
$ws.Range("N116").Value = -4901

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 583.51514
$ws.Range("I3").Value = 591.96875
$ws.Range("J3").Value = 313
$ws.Range("K3").Value = 591.96875
$ws.Range("L3").Value = 313
$ws.Range("M3").Value = -477.96875
$ws.Range("N3").Value = -541
$ws.Range("H20").Value = 3701.182
$ws.Range("I20").Value = 3810.4
$ws.Range("J20").Value = 2609
$ws.Range("K20").Value = 3810.4
$ws.Range("L20").Value = 2609
$ws.Range("M20").Value = -3563.4
$ws.Range("N20").Value = -3103
$ws.Range("H94").Value = 661.63336
$ws.Range("I94").Value = 627.4375
$ws.Range("K94").Value = 627.4375
$ws.Range("M94").Value = -176.4375
$ws.Range("H99").Value = 1137
$ws.Range("I99").Value = 700
$ws.Range("K99").Value = 700
$ws.Range("M99").Value = 798
$ws.Range("H105").Value = 2274582.2
$ws.Range("J105").Value = 2633427
$ws.Range("L105").Value = 2633427
$ws.Range("N105").Value = -2636921

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 31846.666
$ws.Range("I86").Value = 16671.334
$ws.Range("J86").Value = 47022
$ws.Range("K86").Value = 16671.334
$ws.Range("L86").Value = 47022
$ws.Range("M86").Value = -15548.334
$ws.Range("N86").Value = -49268
$ws.Range("H89").Value = 31846.666
$ws.Range("I89").Value = 16671.334
$ws.Range("J89").Value = 47022
$ws.Range("K89").Value = 83356.67
$ws.Range("L89").Value = 235110
$ws.Range("M89").Value = -77740.67
$ws.Range("N89").Value = -246342
$ws.Range("H134").Value = 1100.25
$ws.Range("I134").Value = 971.7143
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2915.1429
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -380.1428999999998
$ws.Range("N134").Value = -11070

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2607.6924
$ws.Range("I59").Value = 450
$ws.Range("K59").Value = 1350
$ws.Range("M59").Value = -810
$ws.Range("H80").Value = 28550.25
$ws.Range("J80").Value = 28550.25
$ws.Range("L80").Value = 85650.75
$ws.Range("N80").Value = -87522.75
$ws.Range("H83").Value = 28550.25
$ws.Range("J83").Value = 28550.25
$ws.Range("L83").Value = 256952.25
$ws.Range("N83").Value = -266312.25
$ws.Range("H122").Value = 625.75
$ws.Range("I122").Value = 442
$ws.Range("J122").Value = 662.5
$ws.Range("K122").Value = 3978
$ws.Range("L122").Value = 5962.5
$ws.Range("M122").Value = -1528
$ws.Range("N122").Value = -10862.5
$ws.Range("H131").Value = 684.2222
$ws.Range("J131").Value = 711.5333000000001
$ws.Range("L131").Value = 2134.5999
$ws.Range("N131").Value = -12214.5999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2394.2
$ws.Range("I102").Value = 2001.1111
$ws.Range("J102").Value = 5932
$ws.Range("K102").Value = 2001.1111
$ws.Range("L102").Value = 5932
$ws.Range("M102").Value = -379.1111000000001
$ws.Range("N102").Value = -9176
$ws.Range("H132").Value = 20991.654
$ws.Range("I132").Value = 3841.348
$ws.Range("J132").Value = 86734.5
$ws.Range("K132").Value = 11524.044
$ws.Range("L132").Value = 260203.5
$ws.Range("M132").Value = -8994.044
$ws.Range("N132").Value = -265263.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1366.7
$ws.Range("I93").Value = 1296.3334
$ws.Range("K93").Value = 1296.3334
$ws.Range("M93").Value = -48.33339999999998
$ws.Range("H100").Value = 2298.1667
$ws.Range("I100").Value = 1763.3334
$ws.Range("J100").Value = 2833
$ws.Range("K100").Value = 1763.3334
$ws.Range("L100").Value = 2833
$ws.Range("M100").Value = -1222.3334
$ws.Range("N100").Value = -3915
$ws.Range("H136").Value = 1637.4783
$ws.Range("I136").Value = 1439.1818
$ws.Range("K136").Value = 4317.5454
$ws.Range("M136").Value = -1767.5454

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1802.4615
$ws.Range("I96").Value = 1712.8
$ws.Range("J96").Value = 2101.3333
$ws.Range("K96").Value = 1712.8
$ws.Range("L96").Value = 2101.3333
$ws.Range("M96").Value = -339.8
$ws.Range("N96").Value = -4847.3333
$ws.Range("H113").Value = 872.65
$ws.Range("I113").Value = 914.3684
$ws.Range("K113").Value = 2743.1052
$ws.Range("M113").Value = -573.1052
